$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1538461538461539
$ws.Cells.Item(2, 3).Value = 0.6470588235294118
$ws.Cells.Item(2, 10).Value = 0.009049773755656109
$ws.Cells.Item(2, 16).Value = 0.1312217194570136
$ws.Cells.Item(2, 19).Value = 0.05882352941176471
$ws.Cells.Item(3, 2).Value = 0.01388888888888889
$ws.Cells.Item(3, 10).Value = 0.02777777777777778
$ws.Cells.Item(3, 16).Value = 0.8125
$ws.Cells.Item(3, 19).Value = 0.1458333333333333
$ws.Cells.Item(4, 10).Value = 0.025
$ws.Cells.Item(4, 16).Value = 0.65
$ws.Cells.Item(4, 19).Value = 0.325
$ws.Cells.Item(6, 2).Value = 0.04661016949152542
$ws.Cells.Item(6, 4).Value = 0.008474576271186441
$ws.Cells.Item(6, 6).Value = 0.06779661016949153
$ws.Cells.Item(6, 10).Value = 0.2457627118644068
$ws.Cells.Item(6, 15).Value = 0.01271186440677966
$ws.Cells.Item(6, 17).Value = 0.1694915254237288
$ws.Cells.Item(6, 18).Value = 0.06779661016949153
$ws.Cells.Item(6, 19).Value = 0.3813559322033898
$ws.Cells.Item(7, 2).Value = 0.05660377358490566
$ws.Cells.Item(7, 4).Value = 0.02830188679245283
$ws.Cells.Item(7, 5).Value = 0.009433962264150943
$ws.Cells.Item(7, 6).Value = 0.06132075471698113
$ws.Cells.Item(7, 10).Value = 0.1462264150943396
$ws.Cells.Item(7, 15).Value = 0.01415094339622642
$ws.Cells.Item(7, 17).Value = 0.2028301886792453
$ws.Cells.Item(7, 18).Value = 0.07547169811320754
$ws.Cells.Item(7, 19).Value = 0.4056603773584906
$ws.Cells.Item(8, 2).Value = 0.0625
$ws.Cells.Item(8, 4).Value = 0.02291666666666667
$ws.Cells.Item(8, 6).Value = 0.05416666666666667
$ws.Cells.Item(8, 10).Value = 0.1041666666666667
$ws.Cells.Item(8, 15).Value = 0.025
$ws.Cells.Item(8, 17).Value = 0.1854166666666667
$ws.Cells.Item(8, 18).Value = 0.08749999999999999
$ws.Cells.Item(8, 19).Value = 0.4583333333333333
$ws.Cells.Item(9, 2).Value = 0.08955223880597014
$ws.Cells.Item(9, 4).Value = 0.007462686567164179
$ws.Cells.Item(9, 6).Value = 0.05970149253731343
$ws.Cells.Item(9, 10).Value = 0.07835820895522388
$ws.Cells.Item(9, 15).Value = 0.01865671641791045
$ws.Cells.Item(9, 17).Value = 0.2574626865671642
$ws.Cells.Item(9, 18).Value = 0.08582089552238806
$ws.Cells.Item(9, 19).Value = 0.4029850746268657
$ws.Cells.Item(10, 2).Value = 0.08502340093603744
$ws.Cells.Item(10, 4).Value = 0.01638065522620905
$ws.Cells.Item(10, 5).Value = 0.0007800312012480499
$ws.Cells.Item(10, 6).Value = 0.0748829953198128
$ws.Cells.Item(10, 10).Value = 0.08970358814352575
$ws.Cells.Item(10, 15).Value = 0.01482059282371295
$ws.Cells.Item(10, 17).Value = 0.2059282371294852
$ws.Cells.Item(10, 18).Value = 0.1060842433697348
$ws.Cells.Item(10, 19).Value = 0.406396255850234
$ws.Cells.Item(11, 7).Value = 0.1682847896440129
$ws.Cells.Item(11, 10).Value = 0.06148867313915857
$ws.Cells.Item(11, 11).Value = 0.2200647249190938
$ws.Cells.Item(11, 12).Value = 0.5469255663430421
$ws.Cells.Item(11, 19).Value = 0.003236245954692557
$ws.Cells.Item(12, 7).Value = 0.7055555555555556
$ws.Cells.Item(12, 10).Value = 0.1888888888888889
$ws.Cells.Item(12, 11).Value = 0.02777777777777778
$ws.Cells.Item(12, 12).Value = 0.03888888888888889
$ws.Cells.Item(12, 19).Value = 0.03888888888888889
$ws.Cells.Item(13, 7).Value = 0.8
$ws.Cells.Item(13, 10).Value = 0.16
$ws.Cells.Item(13, 19).Value = 0.04
$ws.Cells.Item(14, 7).Value = 0.5714285714285714
$ws.Cells.Item(14, 10).Value = 0.2857142857142857
$ws.Cells.Item(14, 19).Value = 0.1428571428571428
$ws.Cells.Item(15, 6).Value = 0.0128755364806867
$ws.Cells.Item(15, 8).Value = 0.1373390557939914
$ws.Cells.Item(15, 9).Value = 0.09012875536480687
$ws.Cells.Item(15, 10).Value = 0.3605150214592275
$ws.Cells.Item(15, 11).Value = 0.05150214592274678
$ws.Cells.Item(15, 13).Value = 0.01716738197424893
$ws.Cells.Item(15, 14).Value = 0.008583690987124463
$ws.Cells.Item(15, 15).Value = 0.04291845493562232
$ws.Cells.Item(15, 19).Value = 0.278969957081545
$ws.Cells.Item(16, 6).Value = 0.01176470588235294
$ws.Cells.Item(16, 8).Value = 0.1941176470588235
$ws.Cells.Item(16, 9).Value = 0.1058823529411765
$ws.Cells.Item(16, 10).Value = 0.3647058823529412
$ws.Cells.Item(16, 11).Value = 0.1235294117647059
$ws.Cells.Item(16, 13).Value = 0.01764705882352941
$ws.Cells.Item(16, 15).Value = 0.02352941176470588
$ws.Cells.Item(16, 19).Value = 0.1588235294117647
$ws.Cells.Item(17, 6).Value = 0.02186878727634195
$ws.Cells.Item(17, 8).Value = 0.1789264413518887
$ws.Cells.Item(17, 9).Value = 0.121272365805169
$ws.Cells.Item(17, 10).Value = 0.3976143141153082
$ws.Cells.Item(17, 11).Value = 0.09145129224652088
$ws.Cells.Item(17, 13).Value = 0.01590457256461232
$ws.Cells.Item(17, 14).Value = 0.005964214711729622
$ws.Cells.Item(17, 15).Value = 0.06163021868787277
$ws.Cells.Item(17, 19).Value = 0.1053677932405567
$ws.Cells.Item(18, 6).Value = 0.01739130434782609
$ws.Cells.Item(18, 8).Value = 0.1826086956521739
$ws.Cells.Item(18, 9).Value = 0.1217391304347826
$ws.Cells.Item(18, 10).Value = 0.4260869565217391
$ws.Cells.Item(18, 11).Value = 0.06521739130434782
$ws.Cells.Item(18, 13).Value = 0.01304347826086956
$ws.Cells.Item(18, 14).Value = 0.004347826086956522
$ws.Cells.Item(18, 15).Value = 0.06521739130434782
$ws.Cells.Item(18, 19).Value = 0.1043478260869565
$ws.Cells.Item(19, 6).Value = 0.01653486700215672
$ws.Cells.Item(19, 8).Value = 0.2070452911574407
$ws.Cells.Item(19, 9).Value = 0.1020848310567937
$ws.Cells.Item(19, 10).Value = 0.3673616103522646
$ws.Cells.Item(19, 11).Value = 0.09992810927390366
$ws.Cells.Item(19, 13).Value = 0.02300503235082674
$ws.Cells.Item(19, 14).Value = 0.0007189072609633358
$ws.Cells.Item(19, 15).Value = 0.07189072609633357
$ws.Cells.Item(19, 19).Value = 0.111430625449317
